# Update all 100 answer cells in the 20x5 table (row-major order)
# matching the document's cell order, using new computed values.
$d = $word.ActiveDocument

$newValues = @(
    "93-19=74",
    "14+12=26",
    "8+8=16",
    "72+16=88",
    "68-65=3",
    "21+11=32",
    "21+24=45",
    "61-45=16",
    "65-5=60",
    "23+27=50",
    "90-1=89",
    "97-89=8",
    "25+25=50",
    "46-42=4",
    "48+42=90",
    "73-25=48",
    "16-12=4",
    "82-21=61",
    "71-52=19",
    "84+12=96",
    "59-24=35",
    "58-3=55",
    "40+2=42",
    "29-15=14",
    "45+8=53",
    "97-57=40",
    "21+72=93",
    "51-8=43",
    "37+33=70",
    "18+40=58",
    "57-5=52",
    "43+56=99",
    "17-1=16",
    "26+8=34",
    "3+94=97",
    "60+5=65",
    "96-94=2",
    "36+28=64",
    "76-43=33",
    "84-35=49",
    "59+1=60",
    "63+23=86",
    "83+4=87",
    "65-18=47",
    "39+54=93",
    "27+66=93",
    "38+50=88",
    "74+13=87",
    "79+11=90",
    "45+27=72",
    "42-39=3",
    "48-21=27",
    "21+38=59",
    "38+0=38",
    "46+36=82",
    "13+37=50",
    "94-92=2",
    "79-28=51",
    "19-6=13",
    "94-47=47",
    "21+65=86",
    "5+18=23",
    "52+38=90",
    "48-37=11",
    "8+86=94",
    "74-51=23",
    "80+2=82",
    "31+61=92",
    "40-15=25",
    "65+23=88",
    "1+27=28",
    "69-21=48",
    "21+15=36",
    "37-9=28",
    "91-34=57",
    "96-64=32",
    "25+5=30",
    "86-1=85",
    "38+30=68",
    "35+26=61",
    "24-1=23",
    "24+69=93",
    "6+0=6",
    "52-23=29",
    "44+5=49",
    "1+83=84",
    "78-10=68",
    "99-14=85",
    "0+98=98",
    "77-51=26",
    "7+75=82",
    "73+7=80",
    "81-15=66",
    "27-26=1",
    "25+45=70",
    "20+74=94",
    "31+24=55",
    "43+24=67",
    "91-70=21",
    "45+9=54"
)

$table = $d.Tables.Item(1)
$rowCount = $table.Rows.Count
$colCount = $table.Columns.Count

$i = 0
for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $table.Cell($r, $c)
        $cell.Range.Text = $newValues[$i]
        $i = $i + 1
    }
}

Write-Output ("Updated " + $i + " cells")
